# Update the "two-digit number divided by one-digit number" worksheet
# table with a freshly generated set of problems/answers. The five
# data rows of the 5-column table (table rows 1, 5, 9, 13, 17 — the
# other rows are blank spacer rows) each get all of their cell text
# replaced in place. Cell text is set directly (rather than via
# Find/Replace, whose scoping to a sub-range is unreliable here) so
# each table cell is targeted unambiguously even though several of
# the old/new problem strings repeat elsewhere in the table.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $table.Cell($row, $col).Range.Text = $newText
}

# Row 1
Set-CellText $tbl 1 1 "50÷9=5, 5"
Set-CellText $tbl 1 2 "62÷5=12, 2"
Set-CellText $tbl 1 3 "16÷6=2, 4"
Set-CellText $tbl 1 4 "24÷6=4, 0"
Set-CellText $tbl 1 5 "71÷7=10, 1"

# Row 5
Set-CellText $tbl 5 1 "24÷8=3, 0"
Set-CellText $tbl 5 2 "68÷8=8, 4"
Set-CellText $tbl 5 3 "13÷5=2, 3"
Set-CellText $tbl 5 4 "95÷9=10, 5"
Set-CellText $tbl 5 5 "72÷5=14, 2"

# Row 9
Set-CellText $tbl 9 1 "61÷2=30, 1"
Set-CellText $tbl 9 2 "73÷6=12, 1"
Set-CellText $tbl 9 3 "76÷7=10, 6"
Set-CellText $tbl 9 4 "10÷5=2, 0"
Set-CellText $tbl 9 5 "42÷3=14, 0"

# Row 13
Set-CellText $tbl 13 1 "83÷3=27, 2"
Set-CellText $tbl 13 2 "75÷5=15, 0"
Set-CellText $tbl 13 3 "82÷9=9, 1"
Set-CellText $tbl 13 4 "78÷9=8, 6"
Set-CellText $tbl 13 5 "91÷8=11, 3"

# Row 17
Set-CellText $tbl 17 1 "38÷8=4, 6"
Set-CellText $tbl 17 2 "33÷3=11, 0"
Set-CellText $tbl 17 3 "37÷7=5, 2"
Set-CellText $tbl 17 4 "84÷8=10, 4"
Set-CellText $tbl 17 5 "58÷3=19, 1"

Write-Host "Done updating table cells."
